# Auto-update draw results: append the 2025-12-06 Pick 3 draw as a new
# row (row 81) at the bottom of the "Results" table, mirroring the
# existing rows' layout (Date, Game, Phase, Result, InsertedAt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every value as literal text (dates, the zero-padded
# "phase" code, and the ISO timestamp are all text, not real
# numbers/dates) so we enter them with a leading apostrophe to stop
# Excel from auto-converting date-/number-looking strings, then strip
# the resulting quote-prefix style so the new cells end up on the
# default "Normal" style, matching every other row in the sheet.
$ws.Range("A81").Value = "'2025-12-06"
$ws.Range("B81").Value = "Pick 3"
$ws.Range("C81").Value = "'251206"
$ws.Range("D81").Value = "7-3-6"
$ws.Range("E81").Value = "'2025-12-06T21:37:10.271+04:00"

$ws.Range("A81:E81").Style = "Normal"
